$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.4375153333333333
$ws.Cells.Item(2, 8).Value = 1.312546
$ws.Cells.Item(2, 9).Value = 0.002535486401940996
$ws.Cells.Item(2, 10).Value = 0.002555908833496712
$ws.Cells.Item(2, 13).Value = 9.363528666666667
$ws.Cells.Item(2, 14).Value = 28.090586
$ws.Cells.Item(2, 15).Value = 0.04175743631338733
$ws.Cells.Item(2, 16).Value = 0.04324026421082073
$ws.Cells.Item(2, 17).Value = 4.096687365772889
$ws.Cells.Item(2, 18).Value = 36.870186291956
$ws.Cells.Item(2, 19).Value = 0.0001058754119525107
$ws.Cells.Item(2, 20).Value = 0.0001105181732591684
$ws.Cells.Item(3, 7).Value = 0.4375153333333333
$ws.Cells.Item(3, 8).Value = 1.312546
$ws.Cells.Item(3, 9).Value = 0.002535486401940996
$ws.Cells.Item(3, 10).Value = 0.002555908833496712
$ws.Cells.Item(3, 15).Value = 0.1749266505387075
$ws.Cells.Item(3, 16).Value = 0.1811383852696593
$ws.Cells.Item(3, 17).Value = 17.16148936492889
$ws.Cells.Item(3, 18).Value = 154.45340428436
$ws.Cells.Item(3, 19).Value = 0.0004435241437779774
$ws.Cells.Item(3, 20).Value = 0.0004629731989960528
$ws.Cells.Item(4, 7).Value = 0.4375153333333333
$ws.Cells.Item(4, 8).Value = 1.312546
$ws.Cells.Item(4, 9).Value = 0.002535486401940996
$ws.Cells.Item(4, 10).Value = 0.002555908833496712
$ws.Cells.Item(4, 13).Value = 71.284935
$ws.Cells.Item(4, 14).Value = 213.854805
$ws.Cells.Item(4, 15).Value = 0.3179011075133629
$ws.Cells.Item(4, 16).Value = 0.3291899382573772
$ws.Cells.Item(4, 17).Value = 31.18825209817
$ws.Cells.Item(4, 18).Value = 280.69426888353
$ws.Cells.Item(4, 19).Value = 0.0008060339352621143
$ws.Cells.Item(4, 20).Value = 0.0008413794710902674
$ws.Cells.Item(5, 7).Value = 0.4375153333333333
$ws.Cells.Item(5, 8).Value = 1.312546
$ws.Cells.Item(5, 9).Value = 0.002535486401940996
$ws.Cells.Item(5, 10).Value = 0.002555908833496712
$ws.Cells.Item(5, 13).Value = 23.0690325
$ws.Cells.Item(5, 14).Value = 46.138065
$ws.Cells.Item(5, 15).Value = 0.1028782726814826
$ws.Cells.Item(5, 16).Value = 0.07102102180339065
$ws.Cells.Item(5, 17).Value = 10.093055443915
$ws.Cells.Item(5, 18).Value = 60.55833266349
$ws.Cells.Item(5, 19).Value = 0.000260846461439077
$ws.Cells.Item(5, 20).Value = 0.0001815232569912487
$ws.Cells.Item(6, 7).Value = 0.4375153333333333
$ws.Cells.Item(6, 8).Value = 1.312546
$ws.Cells.Item(6, 9).Value = 0.002535486401940996
$ws.Cells.Item(6, 10).Value = 0.002555908833496712
$ws.Cells.Item(6, 13).Value = 81.293813
$ws.Cells.Item(6, 14).Value = 243.881439
$ws.Cells.Item(6, 15).Value = 0.3625365329530597
$ws.Cells.Item(6, 16).Value = 0.3754103904587522
$ws.Cells.Item(6, 17).Value = 35.56728969263266
$ws.Cells.Item(6, 18).Value = 320.105607233694
$ws.Cells.Item(6, 19).Value = 0.0009192064495093167
$ws.Cells.Item(6, 20).Value = 0.0009595147331599744
$ws.Cells.Item(7, 9).Value = 0.005820447907772805
$ws.Cells.Item(7, 10).Value = 0.005867329523437988
$ws.Cells.Item(7, 13).Value = 9.363528666666667
$ws.Cells.Item(7, 14).Value = 28.090586
$ws.Cells.Item(7, 15).Value = 0.04175743631338733
$ws.Cells.Item(7, 16).Value = 0.04324026421082073
$ws.Cells.Item(7, 17).Value = 9.40433180341978
$ws.Cells.Item(7, 18).Value = 84.638986230778
$ws.Cells.Item(7, 19).Value = 0.0002430469828242115
$ws.Cells.Item(7, 20).Value = 0.0002537048788054075
$ws.Cells.Item(8, 9).Value = 0.005820447907772805
$ws.Cells.Item(8, 10).Value = 0.005867329523437988
$ws.Cells.Item(8, 15).Value = 0.1749266505387075
$ws.Cells.Item(8, 16).Value = 0.1811383852696593
$ws.Cells.Item(8, 19).Value = 0.001018151457141725
$ws.Cells.Item(8, 20).Value = 0.001062798595720557
$ws.Cells.Item(9, 9).Value = 0.005820447907772805
$ws.Cells.Item(9, 10).Value = 0.005867329523437988
$ws.Cells.Item(9, 13).Value = 71.284935
$ws.Cells.Item(9, 14).Value = 213.854805
$ws.Cells.Item(9, 15).Value = 0.3179011075133629
$ws.Cells.Item(9, 16).Value = 0.3291899382573772
$ws.Cells.Item(9, 17).Value = 71.595570985085
$ws.Cells.Item(9, 18).Value = 644.360138865765
$ws.Cells.Item(9, 19).Value = 0.001850326836104811
$ws.Cells.Item(9, 20).Value = 0.001931465843556238
$ws.Cells.Item(10, 9).Value = 0.005820447907772805
$ws.Cells.Item(10, 10).Value = 0.005867329523437988
$ws.Cells.Item(10, 13).Value = 23.0690325
$ws.Cells.Item(10, 14).Value = 46.138065
$ws.Cells.Item(10, 15).Value = 0.1028782726814826
$ws.Cells.Item(10, 16).Value = 0.07102102180339065
$ws.Cells.Item(10, 17).Value = 23.1695596539575
$ws.Cells.Item(10, 18).Value = 139.017357923745
$ws.Cells.Item(10, 19).Value = 0.0005987976269842155
$ws.Cells.Item(10, 20).Value = 0.000416703738011767
$ws.Cells.Item(11, 9).Value = 0.005820447907772805
$ws.Cells.Item(11, 10).Value = 0.005867329523437988
$ws.Cells.Item(11, 13).Value = 81.293813
$ws.Cells.Item(11, 14).Value = 243.881439
$ws.Cells.Item(11, 15).Value = 0.3625365329530597
$ws.Cells.Item(11, 16).Value = 0.3754103904587522
$ws.Cells.Item(11, 17).Value = 81.64806433911633
$ws.Cells.Item(11, 18).Value = 734.832579052047
$ws.Cells.Item(11, 19).Value = 0.002110125004717843
$ws.Cells.Item(11, 20).Value = 0.00220265646734402
$ws.Cells.Item(12, 7).Value = 99.58055866666666
$ws.Cells.Item(12, 8).Value = 298.741676
$ws.Cells.Item(12, 9).Value = 0.577088694179909
$ws.Cells.Item(12, 10).Value = 0.5817369361698658
$ws.Cells.Item(12, 13).Value = 9.363528666666667
$ws.Cells.Item(12, 14).Value = 28.090586
$ws.Cells.Item(12, 15).Value = 0.04175743631338733
$ws.Cells.Item(12, 16).Value = 0.04324026421082073
$ws.Cells.Item(12, 17).Value = 932.4254157180151
$ws.Cells.Item(12, 18).Value = 8391.828741462135
$ws.Cells.Item(12, 19).Value = 0.02409774439439341
$ws.Cells.Item(12, 20).Value = 0.02515445882117835
$ws.Cells.Item(13, 7).Value = 99.58055866666666
$ws.Cells.Item(13, 8).Value = 298.741676
$ws.Cells.Item(13, 9).Value = 0.577088694179909
$ws.Cells.Item(13, 10).Value = 0.5817369361698658
$ws.Cells.Item(13, 15).Value = 0.1749266505387075
$ws.Cells.Item(13, 16).Value = 0.1811383852696593
$ws.Cells.Item(13, 17).Value = 3906.036127903351
$ws.Cells.Item(13, 18).Value = 35154.32515113016
$ws.Cells.Item(13, 19).Value = 0.100948192336648
$ws.Cells.Item(13, 20).Value = 0.1053748892695283
$ws.Cells.Item(14, 7).Value = 99.58055866666666
$ws.Cells.Item(14, 8).Value = 298.741676
$ws.Cells.Item(14, 9).Value = 0.577088694179909
$ws.Cells.Item(14, 10).Value = 0.5817369361698658
$ws.Cells.Item(14, 13).Value = 71.284935
$ws.Cells.Item(14, 14).Value = 213.854805
$ws.Cells.Item(14, 15).Value = 0.3179011075133629
$ws.Cells.Item(14, 16).Value = 0.3291899382573772
$ws.Cells.Item(14, 17).Value = 7098.59365181702
$ws.Cells.Item(14, 18).Value = 63887.34286635317
$ws.Cells.Item(14, 19).Value = 0.1834571350132335
$ws.Cells.Item(14, 20).Value = 0.1915019460997939
$ws.Cells.Item(15, 7).Value = 99.58055866666666
$ws.Cells.Item(15, 8).Value = 298.741676
$ws.Cells.Item(15, 9).Value = 0.577088694179909
$ws.Cells.Item(15, 10).Value = 0.5817369361698658
$ws.Cells.Item(15, 13).Value = 23.0690325
$ws.Cells.Item(15, 14).Value = 46.138065
$ws.Cells.Item(15, 15).Value = 0.1028782726814826
$ws.Cells.Item(15, 16).Value = 0.07102102180339065
$ws.Cells.Item(15, 17).Value = 2297.22714424949
$ws.Cells.Item(15, 18).Value = 13783.36286549694
$ws.Cells.Item(15, 19).Value = 0.0593698880412414
$ws.Cells.Item(15, 20).Value = 0.04131555162755771
$ws.Cells.Item(16, 7).Value = 99.58055866666666
$ws.Cells.Item(16, 8).Value = 298.741676
$ws.Cells.Item(16, 9).Value = 0.577088694179909
$ws.Cells.Item(16, 10).Value = 0.5817369361698658
$ws.Cells.Item(16, 13).Value = 81.293813
$ws.Cells.Item(16, 14).Value = 243.881439
$ws.Cells.Item(16, 15).Value = 0.3625365329530597
$ws.Cells.Item(16, 16).Value = 0.3754103904587522
$ws.Cells.Item(16, 17).Value = 8095.283314683529
$ws.Cells.Item(16, 18).Value = 72857.54983215177
$ws.Cells.Item(16, 19).Value = 0.2092157343943928
$ws.Cells.Item(16, 20).Value = 0.2183900903518075
$ws.Cells.Item(17, 7).Value = 4.1363315
$ws.Cells.Item(17, 8).Value = 8.272663
$ws.Cells.Item(17, 9).Value = 0.02397084507248554
$ws.Cells.Item(17, 10).Value = 0.01610928107528529
$ws.Cells.Item(17, 13).Value = 9.363528666666667
$ws.Cells.Item(17, 14).Value = 28.090586
$ws.Cells.Item(17, 15).Value = 0.04175743631338733
$ws.Cells.Item(17, 16).Value = 0.04324026421082073
$ws.Cells.Item(17, 17).Value = 38.73065857508634
$ws.Cells.Item(17, 18).Value = 232.383951450518
$ws.Cells.Item(17, 19).Value = 0.00100096103649239
$ws.Cells.Item(17, 20).Value = 0.0006965695699417102
$ws.Cells.Item(18, 7).Value = 4.1363315
$ws.Cells.Item(18, 8).Value = 8.272663
$ws.Cells.Item(18, 9).Value = 0.02397084507248554
$ws.Cells.Item(18, 10).Value = 0.01610928107528529
$ws.Cells.Item(18, 15).Value = 0.1749266505387075
$ws.Cells.Item(18, 16).Value = 0.1811383852696593
$ws.Cells.Item(18, 17).Value = 162.2471343032634
$ws.Cells.Item(18, 18).Value = 973.48280581958
$ws.Cells.Item(18, 19).Value = 0.004193139639112176
$ws.Cells.Item(18, 20).Value = 0.002918009161832258
$ws.Cells.Item(19, 7).Value = 4.1363315
$ws.Cells.Item(19, 8).Value = 8.272663
$ws.Cells.Item(19, 9).Value = 0.02397084507248554
$ws.Cells.Item(19, 10).Value = 0.01610928107528529
$ws.Cells.Item(19, 13).Value = 71.284935
$ws.Cells.Item(19, 14).Value = 213.854805
$ws.Cells.Item(19, 15).Value = 0.3179011075133629
$ws.Cells.Item(19, 16).Value = 0.3291899382573772
$ws.Cells.Item(19, 17).Value = 294.8581221159525
$ws.Cells.Item(19, 18).Value = 1769.148732695715
$ws.Cells.Item(19, 19).Value = 0.007620358196574393
$ws.Cells.Item(19, 20).Value = 0.0053030132425439
$ws.Cells.Item(20, 7).Value = 4.1363315
$ws.Cells.Item(20, 8).Value = 8.272663
$ws.Cells.Item(20, 9).Value = 0.02397084507248554
$ws.Cells.Item(20, 10).Value = 0.01610928107528529
$ws.Cells.Item(20, 13).Value = 23.0690325
$ws.Cells.Item(20, 14).Value = 46.138065
$ws.Cells.Item(20, 15).Value = 0.1028782726814826
$ws.Cells.Item(20, 16).Value = 0.07102102180339065
$ws.Cells.Item(20, 17).Value = 95.42116580427374
$ws.Cells.Item(20, 18).Value = 381.684663217095
$ws.Cells.Item(20, 19).Value = 0.002466079135772741
$ws.Cells.Item(20, 20).Value = 0.001144097602484785
$ws.Cells.Item(21, 7).Value = 4.1363315
$ws.Cells.Item(21, 8).Value = 8.272663
$ws.Cells.Item(21, 9).Value = 0.02397084507248554
$ws.Cells.Item(21, 10).Value = 0.01610928107528529
$ws.Cells.Item(21, 13).Value = 81.293813
$ws.Cells.Item(21, 14).Value = 243.881439
$ws.Cells.Item(21, 15).Value = 0.3625365329530597
$ws.Cells.Item(21, 16).Value = 0.3754103904587522
$ws.Cells.Item(21, 17).Value = 336.2581594670095
$ws.Cells.Item(21, 18).Value = 2017.548956802057
$ws.Cells.Item(21, 19).Value = 0.008690307064533845
$ws.Cells.Item(21, 20).Value = 0.006047591498482639
$ws.Cells.Item(22, 7).Value = 67.39800266666667
$ws.Cells.Item(22, 8).Value = 202.194008
$ws.Cells.Item(22, 9).Value = 0.3905845264378918
$ws.Cells.Item(22, 10).Value = 0.3937305443979143
$ws.Cells.Item(22, 13).Value = 9.363528666666667
$ws.Cells.Item(22, 14).Value = 28.090586
$ws.Cells.Item(22, 15).Value = 0.04175743631338733
$ws.Cells.Item(22, 16).Value = 0.04324026421082073
$ws.Cells.Item(22, 17).Value = 631.0831300454098
$ws.Cells.Item(22, 18).Value = 5679.748170408689
$ws.Cells.Item(22, 19).Value = 0.01630980848772482
$ws.Cells.Item(22, 20).Value = 0.0170250127676361
$ws.Cells.Item(23, 7).Value = 67.39800266666667
$ws.Cells.Item(23, 8).Value = 202.194008
$ws.Cells.Item(23, 9).Value = 0.3905845264378918
$ws.Cells.Item(23, 10).Value = 0.3937305443979143
$ws.Cells.Item(23, 15).Value = 0.1749266505387075
$ws.Cells.Item(23, 16).Value = 0.1811383852696593
$ws.Cells.Item(23, 17).Value = 2643.679016159698
$ws.Cells.Item(23, 18).Value = 23793.11114543728
$ws.Cells.Item(23, 19).Value = 0.06832364296202764
$ws.Cells.Item(23, 20).Value = 0.07131971504358209
$ws.Cells.Item(24, 7).Value = 67.39800266666667
$ws.Cells.Item(24, 8).Value = 202.194008
$ws.Cells.Item(24, 9).Value = 0.3905845264378918
$ws.Cells.Item(24, 10).Value = 0.3937305443979143
$ws.Cells.Item(24, 13).Value = 71.284935
$ws.Cells.Item(24, 14).Value = 213.854805
$ws.Cells.Item(24, 15).Value = 0.3179011075133629
$ws.Cells.Item(24, 16).Value = 0.3291899382573772
$ws.Cells.Item(24, 17).Value = 4804.46223922316
$ws.Cells.Item(24, 18).Value = 43240.16015300844
$ws.Cells.Item(24, 19).Value = 0.1241672535321882
$ws.Cells.Item(24, 20).Value = 0.1296121336003929
$ws.Cells.Item(25, 7).Value = 67.39800266666667
$ws.Cells.Item(25, 8).Value = 202.194008
$ws.Cells.Item(25, 9).Value = 0.3905845264378918
$ws.Cells.Item(25, 10).Value = 0.3937305443979143
$ws.Cells.Item(25, 13).Value = 23.0690325
$ws.Cells.Item(25, 14).Value = 46.138065
$ws.Cells.Item(25, 15).Value = 0.1028782726814826
$ws.Cells.Item(25, 16).Value = 0.07102102180339065
$ws.Cells.Item(25, 17).Value = 1554.80671395242
$ws.Cells.Item(25, 18).Value = 9328.840283714519
$ws.Cells.Item(25, 19).Value = 0.04018266141604518
$ws.Cells.Item(25, 20).Value = 0.02796314557834514
$ws.Cells.Item(26, 7).Value = 67.39800266666667
$ws.Cells.Item(26, 8).Value = 202.194008
$ws.Cells.Item(26, 9).Value = 0.3905845264378918
$ws.Cells.Item(26, 10).Value = 0.3937305443979143
$ws.Cells.Item(26, 13).Value = 81.293813
$ws.Cells.Item(26, 14).Value = 243.881439
$ws.Cells.Item(26, 15).Value = 0.3625365329530597
$ws.Cells.Item(26, 16).Value = 0.3754103904587522
$ws.Cells.Item(26, 17).Value = 5479.040625357502
$ws.Cells.Item(26, 18).Value = 49311.36562821751
$ws.Cells.Item(26, 19).Value = 0.141601160039906
$ws.Cells.Item(26, 20).Value = 0.147810537407958
